# Update "想去人数" (want-to-go count, column F) figures across all four
# sheets to the freshly scraped values (gh-pages rebuild at e2eb17a).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F4" = 240
        "F5" = 1974
        "F7" = 432
        "F8" = 411
        "F10" = 6917
        "F11" = 204
        "F12" = 541
        "F14" = 69
        "F15" = 2393
        "F16" = 1744
        "F17" = 145
        "F18" = 47
        "F19" = 99
        "F21" = 113
        "F22" = 57
        "F24" = 973
        "F25" = 168
        "F26" = 148
        "F27" = 4106
    }
    "演出" = @{
        "F4" = 12
    }
    "本地生活" = @{
        "F2" = 120
        "F3" = 686
    }
    "全部类型" = @{
        "F2" = 120
        "F4" = 686
        "F7" = 240
        "F8" = 1974
        "F10" = 12
        "F12" = 432
        "F13" = 411
        "F15" = 6917
        "F16" = 204
        "F17" = 541
        "F19" = 69
        "F20" = 2393
        "F21" = 1744
        "F22" = 145
        "F23" = 47
        "F24" = 99
        "F26" = 113
        "F27" = 57
        "F29" = 973
        "F30" = 168
        "F31" = 148
        "F32" = 4106
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
